$wb = $excel.ActiveWorkbook

# Workbook revision / view GUIDs (cosmetic, Excel regenerates these on save;
# left to the runtime - only data + selection changes are scripted here).

# --- Sheet "pro": update values in B2:B26 ---
$pro = $wb.Worksheets.Item("pro")
$pro.Range("B2").Value = 990321.512379773
$pro.Range("B3").Value = 1070415.4804183352
$pro.Range("B4").Value = 728301.5619912335
$pro.Range("B5").Value = 834178.0805613386
$pro.Range("B6").Value = 1087372.5359733377
$pro.Range("B7").Value = 762639.0857593106
$pro.Range("B8").Value = 670980.1222569009
$pro.Range("B9").Value = 676210.7756027397
$pro.Range("B10").Value = 1073428.1063382323
$pro.Range("B11").Value = 1170549.5376254548
$pro.Range("B12").Value = 1404638.3305841736
$pro.Range("B13").Value = 1260018.8242749723
$pro.Range("B14").Value = 1326740.5622562943
$pro.Range("B15").Value = 1221894.191260132
$pro.Range("B16").Value = 1070946.9989942736
$pro.Range("B17").Value = 942444.2930038981
$pro.Range("B18").Value = 1063111.2264746665
$pro.Range("B19").Value = 1151682.7509259763
$pro.Range("B20").Value = 1101630.9827816947
$pro.Range("B21").Value = 1275255
$pro.Range("B22").Value = 1098305
$pro.Range("B23").Value = 1098305
$pro.Range("B24").Value = 1237678.4097270947
$pro.Range("B25").Value = 1473204.2116241555
$pro.Range("B26").Value = 1095229.303664732

# --- Sheet "ind": update values in B2:B101 ---
$ind = $wb.Worksheets.Item("ind")
$ind.Range("B2").Value = 1000107.9307687287
$ind.Range("B3").Value = 813274.5810646805
$ind.Range("B4").Value = 626441.2313606321
$ind.Range("B5").Value = 734511.5022678756
$ind.Range("B6").Value = 758323.5958576073
$ind.Range("B7").Value = 840750.0736682169
$ind.Range("B8").Value = 923176.5514788264
$ind.Range("B9").Value = 917681.4529581189
$ind.Range("B10").Value = 486470.23831790633
$ind.Range("B11").Value = 632751.4987911228
$ind.Range("B12").Value = 676976.0659109325
$ind.Range("B13").Value = 544302.3645515034
$ind.Range("B14").Value = 972632.438165192
$ind.Range("B15").Value = 939661.8470409485
$ind.Range("B16").Value = 906691.2559167043
$ind.Range("B17").Value = 879215.7633131678
$ind.Range("B18").Value = 1246685.9071956915
$ind.Range("B19").Value = 1390953.0536615688
$ind.Range("B20").Value = 1050700.3497326132
$ind.Range("B21").Value = 1132360.9986755624
$ind.Range("B22").Value = 631936.3298813395
$ind.Range("B23").Value = 679560.5170608028
$ind.Range("B24").Value = 771145.4924059242
$ind.Range("B25").Value = 879215.763313168
$ind.Range("B26").Value = 756313.3139223305
$ind.Range("B27").Value = 681115.8142589249
$ind.Range("B28").Value = 685454.1315471985
$ind.Range("B29").Value = 482999.32476110605
$ind.Range("B30").Value = 615451.0343192175
$ind.Range("B31").Value = 657580.1229779737
$ind.Range("B32").Value = 664906.9210055833
$ind.Range("B33").Value = 978127.5366858992
$ind.Range("B34").Value = 952483.7435892654
$ind.Range("B35").Value = 923176.5514788261
$ind.Range("B36").Value = 892037.6598614851
$ind.Range("B37").Value = 1089861.2066069476
$ind.Range("B38").Value = 870057.2657786559
$ind.Range("B39").Value = 1245555.6646936545
$ind.Range("B40").Value = 1077039.3100586308
$ind.Range("B41").Value = 1056890.615482704
$ind.Range("B42").Value = 1178515.4627410255
$ind.Range("B43").Value = 1182728.3716069008
$ind.Range("B44").Value = 1005969.3691908162
$ind.Range("B45").Value = 1138767.5834412426
$ind.Range("B46").Value = 920063.7203310019
$ind.Range("B47").Value = 1072679.8052997568
$ind.Range("B48").Value = 1125291.8929049054
$ind.Range("B49").Value = 884627.934831093
$ind.Range("B50").Value = 699749.5955575831
$ind.Range("B51").Value = 1144815.3668506786
$ind.Range("B52").Value = 1105801.8391470877
$ind.Range("B53").Value = 1264249.4610423064
$ind.Range("B54").Value = 968445.6861585483
$ind.Range("B55").Value = 823633.7359168116
$ind.Range("B56").Value = 1176426.1867035949
$ind.Range("B57").Value = 945930.9257733955
$ind.Range("B58").Value = 795098.6167037311
$ind.Range("B59").Value = 916653.2083589279
$ind.Range("B60").Value = 907382.8533483915
$ind.Range("B61").Value = 963249.7034391564
$ind.Range("B62").Value = 431427.5403767492
$ind.Range("B63").Value = 369966.53573376074
$ind.Range("B64").Value = 751449.3271409393
$ind.Range("B65").Value = 746289.713300657
$ind.Range("B66").Value = 720871.0922477156
$ind.Range("B67").Value = 1081164.9894799686
$ind.Range("B68").Value = 922290.2269913473
$ind.Range("B69").Value = 953250.147580778
$ind.Range("B70").Value = 683534.0850874189
$ind.Range("B71").Value = 870506.0079962806
$ind.Range("B72").Value = 1230548.1593332535
$ind.Range("B73").Value = 1199379.9951734564
$ind.Range("B74").Value = 681091.980317427
$ind.Range("B75").Value = 841382.462391895
$ind.Range("B76").Value = 1017309.0362030774
$ind.Range("B77").Value = 926087.6939279702
$ind.Range("B78").Value = 1101113.654620085
$ind.Range("B79").Value = 977314.1192349514
$ind.Range("B80").Value = 1095856.0032444424
$ind.Range("B81").Value = 825716.2229005212
$ind.Range("B82").Value = 797890.9183935059
$ind.Range("B83").Value = 1100285.8594182963
$ind.Range("B84").Value = 1014009.6871968094
$ind.Range("B85").Value = 929522.757850412
$ind.Range("B86").Value = 855870.072861939
$ind.Range("B87").Value = 988135.9287511333
$ind.Range("B88").Value = 1002867.0304461966
$ind.Range("B89").Value = 994836.1907997548
$ind.Range("B90").Value = 1063950.9477360607
$ind.Range("B91").Value = 1069774.433664954
$ind.Range("B92").Value = 1117483.4168141712
$ind.Range("B93").Value = 1071218.9649397545
$ind.Range("B94").Value = 1220653.3189890692
$ind.Range("B95").Value = 1333519.5802183093
$ind.Range("B96").Value = 1315137.49670034
$ind.Range("B97").Value = 1325900.186986321
$ind.Range("B98").Value = 1301124.2243798147
$ind.Range("B99").Value = 1176888.9193337392
$ind.Range("B100").Value = 1215265.4916465275
$ind.Range("B101").Value = 1203895.214489271

# --- Sheet "conso": update values in B2:B26 ---
$conso = $wb.Worksheets.Item("conso")
$conso.Range("B2").Value = 893891.5109372716
$conso.Range("B3").Value = 966185.4092949699
$conso.Range("B4").Value = 657383.7701068305
$conso.Range("B5").Value = 752951.3210933503
$conso.Range("B6").Value = 981491.5448115276
$conso.Range("B7").Value = 688379.0398053506
$conso.Range("B8").Value = 605645.6893311582
$conso.Range("B9").Value = 610366.9176367025
$conso.Range("B10").Value = 968908.2694138418
$conso.Range("B11").Value = 1056573.6405973702
$conso.Range("B12").Value = 1267868.4013958124
$conso.Range("B13").Value = 1137331.2884043553
$conso.Range("B14").Value = 1197555.5984802905
$conso.Range("B15").Value = 1102918.226054572
$conso.Range("B16").Value = 966668.9240178609
$conso.Range("B17").Value = 850678.5493843514
$conso.Range("B18").Value = 959596.1616547599
$conso.Range("B19").Value = 1039543.9330991084
$conso.Range("B20").Value = 994365.6356601053
$conso.Range("B21").Value = 1151084
$conso.Range("B22").Value = 1198807
$conso.Range("B23").Value = 1198807
$conso.Range("B24").Value = 724325.9542113369
$conso.Range("B25").Value = 862162.608595649
$conso.Range("B26").Value = 640960.5307990274

# "VA" (sheet3) is a pure formula sheet (=pro!Bn - conso!Bn); it recalculates
# automatically once pro/conso change - no direct write needed.

# --- Update the selection (active cell) shown on each sheet ---
# Order matters: the final Activate() determines which sheet stays tabSelected.
$ind.Activate()
$ind.Range("D87").Select()

$va = $wb.Worksheets.Item("VA")
$va.Activate()
$va.Range("D87").Select()

$conso.Activate()
$conso.Range("D87").Select()

$pro.Activate()
$pro.Range("D87").Select()

Write-Output "edit complete"
